# V 2.0.2 se arreglo la fechar y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos del paciente
$ws.Range("A6").Value = "ALVIZURES"
$ws.Range("C6").Value = "PÉREZ"
$ws.Range("E6").Value = "MELANY"
$ws.Range("G6").Value = "THAILY"
$ws.Range("I6").Value = "/201762610"

# Dirección actual
$ws.Range("A8").Value = "LT. 1 MZNA. 35 CANTÓN CENTRAL"
$ws.Range("D8").Value = "Z. 24"
$ws.Range("F8").Value = "GUATEMALA"
$ws.Range("H8").Value = "GUATEMALA"

# Fecha de nacimiento, edad, lugar de nacimiento, sexo
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2011-02-08"
$ws.Range("F12").Value = "6 AÑOS"
$ws.Range("H12").Value = "GUATEMALA"
$ws.Range("J12").Value = "Femenino"

# Ocupación, nacionalidad, cédula
$ws.Range("D14").Value = "ESTUDIANTE"
$ws.Range("F14").Value = "GUATEMALTECA"
$ws.Range("H14").Value = "CUI: 2195863140101"

# Nombre del Padre / Madre
$ws.Range("A18").Value = "MÁXIMO DE JESÚS ALVIZURES SAMAYOA"
$ws.Range("F18").Value = "BRENDA CELESTE PÉREZ RODRÍGUEZ"

# Notificar en caso de emergencia
$ws.Range("A20").Value = "MARÍA DEL CARMEN ALVIZURES SAMAYOA"
$ws.Range("F20").Value = "TÍA"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "47057481"

# Fecha y hora de ingreso (reimpresión)
$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "14:20:55"
